# Revert "added ZX correction"
# This reverts a prior edit that had:
#   - replaced "Zephir" with "ZX" in the RSD_Models list (RSD_Models!A4)
#   - changed the Configuration_template RSD model selection (B34) to "Triton"
#   - changed the Configuration_template RSD correction selection (B36) to "ZX"
# Reverting restores "Zephir" in RSD_Models!A4, sets B34 back to "WindCube v2",
# and sets B36 back to "GE".

$wb = $excel.ActiveWorkbook

# --- RSD_Models sheet: restore "Zephir" at A4 (was overwritten with "ZX") ---
$wsModels = $wb.Worksheets.Item("RSD_Models")
$wsModels.Range("A4").Value = "Zephir"
# Restore the original active selection on this sheet.
$wsModels.Activate()
$wsModels.Range("A7").Select()

# --- Configuration_template sheet: restore prior selections ---
$wsConfig = $wb.Worksheets.Item("Configuration_template")
$wsConfig.Range("B34").Value = "WindCube v2"
$wsConfig.Range("B36").Value = "GE"

# Re-activate Configuration_template so it remains the selected/visible tab,
# and restore its original cell selection (B26) and top-left scroll cell.
$wsConfig.Activate()
$wsConfig.Range("B26").Select()
